$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: TN330 / Natalie's - Honey Tangerine / 1 / 14.57 / 14.57
$ws.Range("A20").Value = "TN330"
$ws.Range("B20").Value = "Natalie's - Honey Tangerine"
$ws.Range("C20").Value = "'1"
$ws.Range("D20").Value = "'14.57"
$ws.Range("E20").Value = "'14.57"

# Row 21: TN454 / Natalie's - Orange Mango / 1 / 13.38 / 13.38
$ws.Range("A21").Value = "TN454"
$ws.Range("B21").Value = "Natalie's - Orange Mango"
$ws.Range("C21").Value = "'1"
$ws.Range("D21").Value = "'13.38"
$ws.Range("E21").Value = "'13.38"
